# Workbook was resorted alphabetically (column A) after a new term,
# "Formule", was added between "Dalen" and "Globalegrafiek". Several
# other rows that previously had an empty "answer" column (B) were
# filled in as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row for "Formule" (alphabetically between "Dalen"
#     and "Globalegrafiek", i.e. row 5), shifting everything below down
#     by one row. ---
$ws.Rows("5:5").Insert()

# The insert duplicated the old last row (51) off the bottom of the
# table; drop it again so the sheet keeps its original A1:B50 extent.
$ws.Rows("51:51").Delete()

# --- New row: "Formule" / "Verband, variabelen, tabel verticale as, horizontale as" ---
$ws.Range("A5").Value = "Formule"
$ws.Range("A5").Interior.Color = 65535
$ws.Range("A5").WrapText = $true
$ws.Range("B5").Value = "Verband, variabelen, tabel verticale as, horizontale as"

# --- Fill in answers (column B) for rows that previously had none. Once a
#     term has an answer, its "term" cell (column A) switches from the
#     blank/no-fill look to the yellow highlighted + wrapped look, and the
#     "answer" cell (column B) switches to a yellow, non-wrapped look. ---

# Constant (row 2)
$ws.Range("A2").Interior.Color = 65535
$ws.Range("A2").WrapText = $true
$ws.Range("B2").Value = "Grafiek, horizontaal, formule"
$ws.Range("B2").Interior.Color = 65535

# Dalen (row 4)
$ws.Range("A4").Interior.Color = 65535
$ws.Range("A4").WrapText = $true
$ws.Range("B4").Value = "Grafiek"
$ws.Range("B4").Interior.Color = 65535

# Globalegrafiek (row 6, after the insert)
$ws.Range("A6").Interior.Color = 65535
$ws.Range("A6").WrapText = $true
$ws.Range("B6").Value = "Grafiek, assen"
$ws.Range("B6").Interior.Color = 65535

# Stijgen (row 15, after the insert)
$ws.Range("A15").Interior.Color = 65535
$ws.Range("A15").WrapText = $true
$ws.Range("B15").Value = "Grafiek"
$ws.Range("B15").Interior.Color = 65535

# Tabel (row 16, after the insert) already had an answer; just the answer
# cell's look changes to the yellow, non-wrapped style.
$ws.Range("B16").Interior.Color = 65535

# Vloeiendekromme (row 18, after the insert)
$ws.Range("A18").Interior.Color = 65535
$ws.Range("A18").WrapText = $true
$ws.Range("B18").Value = "Grafiek, punten, liniaal, geodriehoek"
$ws.Range("B18").Interior.Color = 65535

# Woordformule (row 19, after the insert)
$ws.Range("A19").Interior.Color = 65535
$ws.Range("A19").WrapText = $true
$ws.Range("B19").Value = "Verband, variabelen, tabel verticale as, horizontale as"

# --- Update the view state to match: scrolled so row 10 is at the top,
#     with B2 selected. ---
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("B2").Select()

# --- Update the sort range/condition to reflect the new table extent. ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A19:A50"))
$ws.Sort.SetRange($ws.Range("A1:B50"))
$ws.Sort.Header = 0
$ws.Sort.Apply()
